$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 69 (the new weekly data point), shifting the
# existing rows 69-164 down to 70-165. Excel's Insert() copies formatting
# (e.g. the date style on column D) from the row above automatically.
$ws.Rows("69").Insert()

# Populate the newly inserted row 69 with the latest weekly record.
$ws.Cells.Item(69, 1).Value = 9
$ws.Cells.Item(69, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(69, 3).Value = "Metropolitana"
$ws.Cells.Item(69, 4).Value = 44482
$ws.Cells.Item(69, 5).Value = 13
$ws.Cells.Item(69, 6).Value = 300000001
$ws.Cells.Item(69, 7).Value = "Rabanito"
$ws.Cells.Item(69, 8).Value = "Sin especificar"
$ws.Cells.Item(69, 9).Value = "Primera"
$ws.Cells.Item(69, 10).Value = 7900
$ws.Cells.Item(69, 11).Value = 3000
$ws.Cells.Item(69, 12).Value = 4000
$ws.Cells.Item(69, 13).Value = 3494
$ws.Cells.Item(69, 14).Value = "`$/cien unidades (volumen en unidades)"
$ws.Cells.Item(69, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(69, 16).Value = 35
$ws.Cells.Item(69, 17).Value = 100
$ws.Cells.Item(69, 18).Value = "Hortaliza"
